# Update countries & provincias Spain
# - Refresh case counts for a subset of countries (new pull of data).
# - Update the "last updated" timestamp.
# - Re-sort the country table by "Casos totales" (column B) descending,
#   which is what naturally moves Bosnia y Herzegovina, Croacia and
#   Montenegro ahead of their neighbours after their counts increased.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last data pull timestamp (row 1).
$ws.Range("A1").Value = "Datos actualizados a 31 de Julio de 2020 a las 15:09"

# New stats per country: Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes
$updates = @{
    "Estados Unidos"       = @(4635886, 901, 2286277, 2194279, 0, 45, 155330)
    "India"                = @(1663174, 23824, 1069501, 557693, 0, 194, 35980)
    "Alemania"             = @(209980, 327, 192000, 8759, 0, 0, 9221)
    "Paises Bajos"         = @(54301, 338, 0, 0, 0, 0, 6147)
    "Azerbaiyan"           = @(31878, 318, 25882, 5548, 0, 7, 448)
    "Uzbekistan"           = @(23773, 502, 14204, 9430, 0, 3, 139)
    "Dinamarca"            = @(13789, 64, 12578, 596, 0, 0, 615)
    "Bosnia y Herzegovina" = @(11876, 432, 5959, 5578, 0, 11, 339)
    "Zambia"               = @(5963, 408, 3803, 2009, 0, 2, 151)
    "Croacia"              = @(5139, 68, 4267, 727, 0, 1, 145)
    "Montenegro"           = @(3073, 57, 1005, 2020, 0, 1, 48)
    "Islandia"             = @(1885, 13, 1825, 50, 0, 0, 10)
}

foreach ($country in $updates.Keys) {
    # LookAt:=xlWhole (1) so e.g. "Estados Unidos" doesn't match the
    # substring inside "Islas Virgenes de los Estados Unidos".
    $found = $ws.Range("A4:A219").Find($country, [Type]::Missing, [Type]::Missing, 1)
    $row = $found.Row
    $vals = $updates[$country]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
    $ws.Cells.Item($row, 8).Value = $vals[6]
}

# Re-sort the whole country table by Casos totales (column B) descending.
$sortRange = $ws.Range("A4:H219")
$keyRange = $ws.Range("B4:B219")
$sortRange.Sort($keyRange, 2)
